$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 4 - "Using Oldest Entry Date" case for Renewal
$ws.Range("A4").Value = "AAA_CSA"
$ws.Range("B4").Value = "CHOICE"
$ws.Range("C4").Value = "CA"
$ws.Range("D4").Value = 20000102
$ws.Range("E4").Value = 20300102
$ws.Range("F4").Value = "SYMBOL_2000_ENTRY_DATE"

# Copy styling from row 3 to row 4 to match formatting
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)
$ws.Range("D3:E3").Copy()
$ws.Range("D4:E4").PasteSpecial(-4122)
$ws.Range("F3").Copy()
$ws.Range("F4").PasteSpecial(-4122)

$ws.Range("F4").Select()
